$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Damage table tweaks (rows 1-6) ---
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 7
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 10

# Halberd damage row gains values
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 4.5
$ws.Range("D6").Value = 12

# --- Insert a blank separator row before the "Time" header block ---
$ws.Rows.Item(7).Insert()

# --- Time table: Halberd row (now row 13) gains values ---
$ws.Range("B13").Value = 0.4
$ws.Range("C13").Value = 0.4
$ws.Range("D13").Value = 1

# --- Insert a blank separator row before the "DPS" header block ---
$ws.Rows.Item(14).Insert()

# --- DPS table formulas (now rows 17-20), one shared relative formula ---
$ws.Range("B17:D20").FormulaR1C1 = "=PRODUCT(R[-14]C,1/R[-7]C)"

# --- Column A width (auto-fit to content, matches the widened "GreatSword" label) ---
$ws.Columns.Item(1).AutoFit()

# --- Selection as left by the editor ---
$ws.Range("B14").Select()
